# ToDo.xlsx update: mark the Python-model rows as "Done" and layer on
# conditional formatting that highlights "To-Do" (red) / "Done" (green)
# cells in the status columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Mark Python-model rows as Done -----------------------------------
# B6 ("Perceptron Forward Python") was hand-marked Done with a direct
# green (Accent 6 theme) fill.
$ws.Range("B6").Value = "Done"
$ws.Range("B6").Interior.ThemeColor = 10

# The remaining Python-model rows just get the text swapped to "Done",
# keeping their existing (red) fill style.
$pythonDoneRows = @(8, 10, 14, 16, 18, 22, 24, 26)
foreach ($r in $pythonDoneRows) {
    $ws.Cells.Item($r, 2).Value = "Done"
}

# --- Conditional formatting --------------------------------------------
# Whole status-grid rule set: "To-Do" -> red, "Done" -> green.
$rngAll = $ws.Range("B3:F29")
$fcsAll = $rngAll.FormatConditions

$fTodoAll = $fcsAll.Add(9, 0, [Type]::Missing, [Type]::Missing, "To-Do")
$fTodoAll.Font.Color = 393372
$fTodoAll.Interior.Color = 13551615

$fDoneAll = $fcsAll.Add(9, 0, [Type]::Missing, [Type]::Missing, "Done")
$fDoneAll.Font.Color = 24832
$fDoneAll.Interior.Color = 13561798
$fDoneAll.SetFirstPriority()

# Earlier experimentation on B7 (added/removed a few times while testing
# the rule) before settling on a plain "To-Do" rule there too.
$rngB7 = $ws.Range("B7")
$fcsB7 = $rngB7.FormatConditions

$tmp1 = $fcsB7.Add(9, 0, [Type]::Missing, [Type]::Missing, "To-Do")
$tmp1.Font.Color = 393372
$tmp1.Interior.Color = 13551615
$fcsB7.Delete()

$tmp2 = $fcsB7.Add(9, 0, [Type]::Missing, [Type]::Missing, "To-Do")
$tmp2.Font.Color = 393372
$tmp2.Interior.Color = 13551615
$fcsB7.Delete()

$fB7 = $fcsB7.Add(9, 0, [Type]::Missing, [Type]::Missing, "To-Do")
$fB7.Font.Color = 393372
$fB7.Interior.Color = 13551615

# One more stray rule created and removed elsewhere, leaving an unused
# dxf behind (matches the accumulated dxf count in the saved file).
$rngTemp = $ws.Range("A1")
$fcsTemp = $rngTemp.FormatConditions
$tmp3 = $fcsTemp.Add(9, 0, [Type]::Missing, [Type]::Missing, "To-Do")
$tmp3.Font.Color = 393372
$tmp3.Interior.Color = 13551615
$fcsTemp.Delete()

# --- Final selection -----------------------------------------------------
$null = $ws.Range("M20").Select()
